# Helper: assign a horizontal array of values to a single-row range using a 2D COM-safe array
function Set-RowValues {
    param($ws, $rangeAddr, $values)
    $n = $values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) { $arr[0,$i] = $values[$i] }
    $ws.Range($rangeAddr).Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Reuse existing cell styles (avoid creating new style entries) by copying format
# from the (now-shifted) first data column of each row-style-group into the two new columns
$ws.Range("F7").Copy() | Out-Null
$ws.Range("D7:E7").PasteSpecial(-4122) | Out-Null
$ws.Range("D38:E38").PasteSpecial(-4122) | Out-Null
$ws.Range("D80:E80").PasteSpecial(-4122) | Out-Null

$ws.Range("F8").Copy() | Out-Null
$ws.Range("D8:E35").PasteSpecial(-4122) | Out-Null
$ws.Range("D39:E77").PasteSpecial(-4122) | Out-Null
$ws.Range("D81:E102").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Set final data values for D:M across all data rows
Set-RowValues $ws "D7:M7" @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowValues $ws "D8:M8" @(2011400,2067200,2104000,1993900,1996600,2000200,3682200,1772600,1793300,1855200)
Set-RowValues $ws "D9:M9" @(1621100,1670000,1704100,1614000,1613500,1614800,2966100,1422600,1444700,1490200)
Set-RowValues $ws "D10:M10" @(390300,397200,399900,379900,383100,385400,716100,350000,348600,365000)
Set-RowValues $ws "D12:M12" @("NA","NA","NA","NA","NA","NA","NA","NA","NA","NA")
Set-RowValues $ws "D13:M13" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D14:M14" @(0,"NA","NA","NA",0,0,0,"NA",0,123900)
Set-RowValues $ws "D15:M15" @(15700,15600,15800,15900,16300,16100,31700,16000,16600,16900)
Set-RowValues $ws "D17:M17" @(1920900,1969700,2012800,1920700,1915700,1911400,3532900,1706000,1711300,1886600)
Set-RowValues $ws "D18:M18" @(90500,97500,91200,73200,80900,88800,149300,66600,82000,-31400)
Set-RowValues $ws "D20:M20" @(-16800,-17000,-17700,-19800,-17100,-16900,-32600,-16300,-17500,-20800)
Set-RowValues $ws "D21:M21" @(89300,96100,89300,69300,80100,88000,148300,66300,81100,-35200)
Set-RowValues $ws "D22:M22" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D23:M23" @(73700,80500,73400,53500,63800,71900,116700,50400,64500,-52200)
Set-RowValues $ws "D24:M24" @(19900,13800,15800,10500,15200,18400,29300,12600,16800,-21100)
Set-RowValues $ws "D25:M25" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D26:M26" @(53800,66600,57700,43000,48600,53600,87300,37800,47800,-31000)
Set-RowValues $ws "D27:M27" @(53800,66800,57900,44400,48900,53700,87200,37700,47300,-31600)
Set-RowValues $ws "D28:M28" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D29:M29" @(4300,"NA","NA","NA",-26400,"NA","NA","NA","NA","NA")
Set-RowValues $ws "D30:M30" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D31:M31" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D32:M32" @(16800,17000,17700,19800,17100,16900,32600,16300,17500,20800)
Set-RowValues $ws "D33:M33" @(58100,66800,57900,44400,22500,53700,87200,37700,47300,-31600)
Set-RowValues $ws "D34:M34" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D35:M35" @(58100,66800,57900,44400,22500,53700,87200,37700,47300,-31600)
Set-RowValues $ws "D38:M38" @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowValues $ws "D41:M41" @(96300,142800,110900,123900,118000,94100,87800,103000,110100,112800)
Set-RowValues $ws "D42:M42" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D43:M43" @(1288500,1346600,1319700,1276400,1294600,1319200,1214300,1127800,1125100,1170600)
Set-RowValues $ws "D44:M44" @(948700,926800,935200,949500,956100,925000,866300,850100,821400,832500)
Set-RowValues $ws "D45:M45" @(52100,90400,81100,68900,40200,71800,137100,122400,46400,127800)
Set-RowValues $ws "D46:M46" @(2385600,2506500,2446900,2418700,2408800,2410100,2305500,2203400,2103100,2243800)
Set-RowValues $ws "D47:M47" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D48:M48" @(160900,157100,157500,156400,156400,157600,155200,156600,157600,160900)
Set-RowValues $ws "D49:M49" @(2038600,2089300,2085300,2109000,2139000,2155100,2122000,2110900,2124300,2143100)
Set-RowValues $ws "D50:M50" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D51:M51" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D52:M52" @(19900,25800,25300,29300,31200,41700,41000,40500,46800,61500)
Set-RowValues $ws "D53:M53" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D54:M54" @(4605000,4778800,4715000,4713300,4735500,4764500,4623700,4511200,4431800,4609400)
Set-RowValues $ws "D57:M57" @(794300,813400,818200,805400,799500,838400,769000,712400,684700,699500)
Set-RowValues $ws "D58:M58" @(74000,62600,65000,70000,72900,67100,54700,62700,51500,72600)
Set-RowValues $ws "D59:M59" @(193600,175600,142700,146200,168500,138000,140700,144900,137600,187900)
Set-RowValues $ws "D60:M60" @(1061900,1051500,1025900,1021600,1041000,1043500,964400,920100,873800,960000)
Set-RowValues $ws "D61:M61" @(1167300,1229300,1261700,1292100,1313300,1368300,1334500,1309800,1363100,1418700)
Set-RowValues $ws "D62:M62" @(246100,266100,261600,266700,265100,245400,229600,223900,231300,238000)
Set-RowValues $ws "D63:M63" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D64:M64" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D65:M65" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D66:M66" @(2469700,2541400,2543900,2575400,2615700,2654000,2525400,2450600,2464900,2613000)
Set-RowValues $ws "D68:M68" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D69:M69" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D70:M70" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D71:M71" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D72:M72" @(2307500,2249300,2182500,2124500,2079700,2057000,2044700,1994600,1914800,1909200)
Set-RowValues $ws "D73:M73" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D74:M74" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D75:M75" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D76:M76" @(2135300,2237400,2171100,2138000,2119700,2110500,2098300,2060700,1966900,1996300)
Set-RowValues $ws "D77:M77" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D80:M80" @(43465,43373,43281,43190,43100,43008,42916,42825,42735,42643)
Set-RowValues $ws "D81:M81" @(58100,66800,57900,44400,22500,53700,87200,37700,47300,-31600)
Set-RowValues $ws "D83:M83" @(15700,15600,15800,15900,16300,16100,31700,16000,16600,16900)
Set-RowValues $ws "D84:M84" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D85:M85" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D86:M86" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D87:M87" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D88:M88" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D89:M89" @(122300,87600,33800,53000,68000,14300,66800,47600,83000,78600)
Set-RowValues $ws "D91:M91" @(-12500,-7400,-8700,-7700,-5500,-6200,-9800,-4500,-4800,-6100)
Set-RowValues $ws "D92:M92" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D93:M93" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D94:M94" @(-14000,4900,-8600,-16400,7200,-6200,-6300,-4500,-2700,-1700)
Set-RowValues $ws "D96:M96" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D97:M97" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D98:M98" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D99:M99" @(0,0,0,0,0,0,0,0,0,0)
Set-RowValues $ws "D100:M100" @(-150300,-61100,-34900,-28800,-49000,-5600,-86500,-50700,-78200,-121900)
Set-RowValues $ws "D101:M101" @(-4400,300,-3200,-1800,-2300,3700,3800,400,-4900,-2400)
Set-RowValues $ws "D102:M102" @(-46400,31800,-13000,5900,23900,6300,-22300,-7100,-2700,-47500)
